$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.561.73'
$ws.Range("E2").Value = '  -4.47%  '
$ws.Range("D3").Value = '3.777.22'
$ws.Range("E3").Value = '  -5.33%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.86'
$ws.Range("E5").Value = '  -1.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.34'
$ws.Range("E6").Value = '  +1.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.653'
$ws.Range("E7").Value = '  -4.57%  '
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.728'
$ws.Range("E9").Value = '  -2.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("E10").Value = '  +1.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '51.55'
$ws.Range("E11").Value = '  -4.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000314'
$ws.Range("E12").Value = '  -1.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.92'
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("D14").Value = '4.402.52'
$ws.Range("E14").Value = '  -4.90%  '
$ws.Range("D15").Value = '3.808.78'
$ws.Range("E15").Value = '  -4.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.36'
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.18'
$ws.Range("E17").Value = '  -7.44%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.45'
$ws.Range("E18").Value = '  -4.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.128'
$ws.Range("D20").Value = '69.657.23'
$ws.Range("E20").Value = '  -4.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '430.11'
$ws.Range("E21").Value = '  -0.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.61'
$ws.Range("E22").Value = '  -3.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '91.75'
$ws.Range("E23").Value = '  -4.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.20'
$ws.Range("E24").Value = '  -6.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.55'
$ws.Range("E25").Value = '  -4.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.03'
$ws.Range("E26").Value = '  -2.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.88'
$ws.Range("E27").Value = '  -12.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.93'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.17'
$ws.Range("E29").Value = '  -2.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.28'
$ws.Range("E30").Value = '  -5.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.84'
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.17'
$ws.Range("E32").Value = '  -4.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '47.07'
$ws.Range("E33").Value = '  -3.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.122'
$ws.Range("E34").Value = '  -6.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '68.23'
$ws.Range("E35").Value = '  -3.37%  '
$ws.Range("D36").Value = '0.0₃0959'
$ws.Range("E36").Value = '  +9.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '621.00'
$ws.Range("E37").Value = '  -7.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.415'
$ws.Range("E38").Value = '  -5.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.142'
$ws.Range("E41").Value = '  -3.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.17'
$ws.Range("E42").Value = '  -5.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.06'
$ws.Range("E43").Value = '  +16.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0459'
$ws.Range("E44").Value = '  -6.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.73'
$ws.Range("E45").Value = '  +4.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.69'
$ws.Range("E46").Value = '  -9.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.141'
$ws.Range("E47").Value = '  -6.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.80'
$ws.Range("E48").Value = '  -15.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.23'
$ws.Range("E49").Value = '  -6.28%  '
$ws.Range("D50").Value = '2.812.23'
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000265'
$ws.Range("E51").Value = '  -0.95%  '
